# Update the "Datasets Used" worksheet to include the final lung/muscle
# datasets (pre Seurat V5 migration).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the header for column C to reflect the additional accession sources
# (Dryad / SingleCellPortal / HCA) now referenced by the new rows below.
$ws.Range("C2").Value = "GEO/EGA Accession/Dryad/SingleCellPortal/HCAProject"

# New Muscle dataset rows (21-23).
$ws.Range("A21").Value = "Human"
$ws.Range("B21").Value = "Muscle"
$ws.Range("C21").Value = "SCP1479"
$ws.Range("D21").Value = "Eraslan/Regev"
$ws.Range("E21").Value = "SN"
$ws.Range("F21").Value = "209k"
$ws.Range("G21").Value = "16 donors"

$ws.Range("A22").Value = "Mouse"
$ws.Range("B22").Value = "Muscle"
$ws.Range("C22").Value = "doi:10.5061/dryad.t4b8gtj34"
$ws.Range("D22").Value = "McKellar/Cosgrove"
$ws.Range("E22").Value = "SC+ SN"
$ws.Range("F22").Value = "365k"
$ws.Range("G22").Value = "Is integrated from 20 sc/sn datasets across 111 samples,"

$ws.Range("A23").Value = "Human"
$ws.Range("B23").Value = "Muscle"
$ws.Range("C23").Value = "10201832-7c73-4033-9b65-3ef13d81656a"
$ws.Range("D23").Value = "Quake"
$ws.Range("E23").Value = "SC"
$ws.Range("F23").Value = "31k"
$ws.Range("G23").Value = "5 donors"

# Match the author's final selection position after entering the new rows.
$ws.Range("C24").Select()
